$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.501.21'
$ws.Range("E2").Value = '  +0.91%  '

$ws.Range("D3").Value = '1.873.41'
$ws.Range("E3").Value = '  +0.75%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.72%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.62%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5083'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.36%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3900'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.32%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08356'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.13%  '

$ws.Range("E10").Value = '  -0.58%  '

$ws.Range("E11").Value = '  -0.49%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.223'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.42%  '

$ws.Range("D13").Value = '1.866.95'
$ws.Range("E13").Value = '  +0.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.78%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.264'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.12%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.008'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.71%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.19'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06726'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.34%  '

$ws.Range("E20").Value = '  +0.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.007'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.70%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.911'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.26%  '

$ws.Range("D23").Value = '28.515.91'
$ws.Range("E23").Value = '  +0.83%  '

$ws.Range("E24").Value = '  +0.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.232'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.80%  '

$ws.Range("D26").Value = '2.086.87'
$ws.Range("E26").Value = '  +0.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.400'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1043'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.23%  '

$ws.Range("E32").Value = '  +1.90%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.756'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.614'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.55%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02454'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.56%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06530'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2159'
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.807'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.030'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.189'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.91%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.241'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6383'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("E44").Value = '  -0.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5994'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.686'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.20%  '

$ws.Range("E48").Value = '  +2.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.213'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '121.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.145'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -10.80%  '

